# Update leve crafting-profit market data cells (scheduled market-price refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 451.83334
$ws.Range("I38").Value = 177.75
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 533.25
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -161.25
$ws.Range("N38").Value = -3744
$ws.Range("H40").Value = 1962.6666
$ws.Range("I40").Value = 1944
$ws.Range("K40").Value = 1944
$ws.Range("M40").Value = -1769
$ws.Range("H61").Value = 221
$ws.Range("I61").Value = 263.75
$ws.Range("J61").Value = 50
$ws.Range("K61").Value = 791.25
$ws.Range("L61").Value = 150
$ws.Range("M61").Value = -619.25
$ws.Range("N61").Value = -494
$ws.Range("H70").Value = 2285.8333
$ws.Range("I70").Value = 1930
$ws.Range("J70").Value = 2357
$ws.Range("K70").Value = 5790
$ws.Range("L70").Value = 7071
$ws.Range("M70").Value = -5520
$ws.Range("N70").Value = -7611
$ws.Range("H73").Value = 2285.8333
$ws.Range("I73").Value = 1930
$ws.Range("J73").Value = 2357
$ws.Range("K73").Value = 5790
$ws.Range("L73").Value = 7071
$ws.Range("M73").Value = -4854
$ws.Range("N73").Value = -8943
$ws.Range("H76").Value = 9070
$ws.Range("I76").Value = 14188.889
$ws.Range("J76").Value = 3311.25
$ws.Range("K76").Value = 14188.889
$ws.Range("L76").Value = 3311.25
$ws.Range("M76").Value = -13873.889
$ws.Range("N76").Value = -3941.25
$ws.Range("H79").Value = 9070
$ws.Range("I79").Value = 14188.889
$ws.Range("J79").Value = 3311.25
$ws.Range("K79").Value = 14188.889
$ws.Range("L79").Value = 3311.25
$ws.Range("M79").Value = -13096.889
$ws.Range("N79").Value = -5495.25
$ws.Range("H97").Value = 1053
$ws.Range("J97").Value = 1053
$ws.Range("L97").Value = 3159
$ws.Range("N97").Value = -4151
$ws.Range("H141").Value = 2224.875
$ws.Range("I141").Value = 1550
$ws.Range("K141").Value = 4650
$ws.Range("M141").Value = 530

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 6825.8125
$ws.Range("I88").Value = 2996.5
$ws.Range("J88").Value = 9123.4
$ws.Range("K88").Value = 2996.5
$ws.Range("L88").Value = 9123.4
$ws.Range("M88").Value = -2590.5
$ws.Range("N88").Value = -9935.4
$ws.Range("H91").Value = 6825.8125
$ws.Range("I91").Value = 2996.5
$ws.Range("J91").Value = 9123.4
$ws.Range("K91").Value = 2996.5
$ws.Range("L91").Value = 9123.4
$ws.Range("M91").Value = -1592.5
$ws.Range("N91").Value = -11931.4
$ws.Range("H122").Value = 5052331
$ws.Range("I122").Value = 1896.2941
$ws.Range("J122").Value = 22223808
$ws.Range("K122").Value = 5688.8823
$ws.Range("L122").Value = 66671424
$ws.Range("M122").Value = -3238.8823
$ws.Range("N122").Value = -66676324

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2346.889
$ws.Range("I70").Value = 1105.5
$ws.Range("J70").Value = 3340
$ws.Range("K70").Value = 3316.5
$ws.Range("L70").Value = 10020
$ws.Range("M70").Value = -3001.5
$ws.Range("N70").Value = -10650
$ws.Range("H73").Value = 2346.889
$ws.Range("I73").Value = 1105.5
$ws.Range("J73").Value = 3340
$ws.Range("K73").Value = 3316.5
$ws.Range("L73").Value = 10020
$ws.Range("M73").Value = -2224.5
$ws.Range("N73").Value = -12204
$ws.Range("H100").Value = 3023.3333
$ws.Range("J100").Value = 3124.5
$ws.Range("L100").Value = 9373.5
$ws.Range("N100").Value = -10995.5
$ws.Range("H106").Value = 3286.6667
$ws.Range("J106").Value = 3286.6667
$ws.Range("L106").Value = 9860.000100000001
$ws.Range("N106").Value = -11752.0001
$ws.Range("H112").Value = 13336258
$ws.Range("I112").Value = 2118.1428
$ws.Range("J112").Value = 18521758
$ws.Range("K112").Value = 6354.428400000001
$ws.Range("L112").Value = 55565274
$ws.Range("M112").Value = -5246.428400000001
$ws.Range("N112").Value = -55567490
$ws.Range("H123").Value = 966.6667
$ws.Range("I123").Value = 966.6667
$ws.Range("K123").Value = 2900.0001
$ws.Range("M123").Value = -450.0001000000002
$ws.Range("H131").Value = 1394.075
$ws.Range("I131").Value = 639.2308
$ws.Range("J131").Value = 1757.5186
$ws.Range("K131").Value = 1917.6924
$ws.Range("L131").Value = 5272.5558
$ws.Range("M131").Value = 3122.3076
$ws.Range("N131").Value = -15352.5558

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7320.846
$ws.Range("I80").Value = 15187.5
$ws.Range("J80").Value = 3824.5557
$ws.Range("K80").Value = 15187.5
$ws.Range("L80").Value = 3824.5557
$ws.Range("M80").Value = -14189.5
$ws.Range("N80").Value = -5820.5557
$ws.Range("H83").Value = 7320.846
$ws.Range("I83").Value = 15187.5
$ws.Range("J83").Value = 3824.5557
$ws.Range("K83").Value = 75937.5
$ws.Range("L83").Value = 19122.7785
$ws.Range("M83").Value = -70945.5
$ws.Range("N83").Value = -29106.7785
$ws.Range("H97").Value = 1123.8572
$ws.Range("I97").Value = 1242.6666
$ws.Range("J97").Value = 826.8333
$ws.Range("K97").Value = 1242.6666
$ws.Range("L97").Value = 826.8333
$ws.Range("M97").Value = -746.6666
$ws.Range("N97").Value = -1818.8333
$ws.Range("H122").Value = 3563.4285
$ws.Range("I122").Value = 2199.8333
$ws.Range("J122").Value = 4586.125
$ws.Range("K122").Value = 6599.499899999999
$ws.Range("L122").Value = 13758.375
$ws.Range("M122").Value = -4149.499899999999
$ws.Range("N122").Value = -18658.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3667.8235
$ws.Range("I16").Value = 1321.0834
$ws.Range("J16").Value = 9300
$ws.Range("K16").Value = 1321.0834
$ws.Range("L16").Value = 9300
$ws.Range("M16").Value = -1151.0834
$ws.Range("N16").Value = -9640
$ws.Range("H82").Value = 3918.625
$ws.Range("I82").Value = 1558
$ws.Range("J82").Value = 5335
$ws.Range("K82").Value = 1558
$ws.Range("L82").Value = 5335
$ws.Range("M82").Value = -1197
$ws.Range("N82").Value = -6057
$ws.Range("H85").Value = 3918.625
$ws.Range("I85").Value = 1558
$ws.Range("J85").Value = 5335
$ws.Range("K85").Value = 1558
$ws.Range("L85").Value = 5335
$ws.Range("M85").Value = -310
$ws.Range("N85").Value = -7831

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3429.3635
$ws.Range("I62").Value = 3225
$ws.Range("J62").Value = 3674.6
$ws.Range("K62").Value = 3225
$ws.Range("L62").Value = 3674.6
$ws.Range("M62").Value = -2601
$ws.Range("N62").Value = -4922.6
$ws.Range("H65").Value = 3429.3635
$ws.Range("I65").Value = 3225
$ws.Range("J65").Value = 3674.6
$ws.Range("K65").Value = 16125
$ws.Range("L65").Value = 18373
$ws.Range("M65").Value = -13005
$ws.Range("N65").Value = -24613

Write-Host "Applied 178 cell updates across ALC, ARM, CUL, GSM, LTW, WVR"
